# Update the "dSF" (column F) values for the irvin_cole 2022 data sheet.
# These rows correspond to the "push all data, mean calculation / repull data"
# commit, which recomputed the dSF column (previously mirroring dS0 / 0-filled)
# with the actual recalculated deltas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    4  = -7
    5  = 3
    6  = 5
    7  = -1
    8  = 5
    9  = 2
    10 = -4
    11 = -5
    12 = -1
    13 = -2
    14 = 1
    15 = 3
    16 = -1
    17 = -2
    18 = 4
    19 = -2
    20 = -1
    22 = -3
    23 = 0
    24 = -1
    25 = 3
    26 = -3
    28 = 2
    29 = -1
    30 = 5
    31 = -4
    32 = -1
    33 = 1
    34 = -3
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
